$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-31: columns D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg) are shuffled
# between rows (a full-row permutation of the weekly price records).
$rows = @{
    2 = @{ D = 44379; L = "Primera"; M = 150; N = 700; O = 800; P = 747; S = 747 }
    3 = @{ D = 44379; L = "Segunda"; M = 140; N = 500; O = 600; P = 543; S = 543 }
    4 = @{ D = 44389; L = "Primera"; M = 140; N = 750; O = 800; P = 775; S = 775 }
    5 = @{ D = 44389; L = "Segunda"; M = 120; N = 600; O = 700; P = 650; S = 650 }
    6 = @{ D = 44403; L = "Primera"; M = 100; N = 1200; O = 1300; P = 1250; S = 1250 }
    7 = @{ D = 44403; L = "Segunda"; M = 120; N = 950; O = 1000; P = 975; S = 975 }
    8 = @{ D = 44372; L = "Primera"; M = 900; N = 750; O = 800; P = 772; S = 772 }
    9 = @{ D = 44372; L = "Segunda"; M = 900; N = 600; O = 650; P = 628; S = 628 }
    10 = @{ D = 44326; L = "Primera"; M = 160; N = 600; O = 700; P = 650; S = 650 }
    11 = @{ D = 44344; L = "Primera"; M = 140; N = 1000; O = 1200; P = 1100; S = 1100 }
    12 = @{ D = 44344; L = "Segunda"; M = 120; N = 800; O = 850; P = 825; S = 825 }
    13 = @{ D = 44316; L = "Primera"; M = 140; N = 1100; O = 1200; P = 1150; S = 1150 }
    14 = @{ D = 44330; L = "Primera"; M = 200; N = 1200; O = 1300; P = 1250; S = 1250 }
    15 = @{ D = 44330; L = "Segunda"; M = 100; N = 1000; O = 1100; P = 1050; S = 1050 }
    16 = @{ D = 44386; L = "Primera"; M = 160; N = 700; O = 750; P = 725; S = 725 }
    17 = @{ D = 44386; L = "Segunda"; M = 200; N = 600; O = 650; P = 625; S = 625 }
    18 = @{ D = 44309; L = "Primera"; M = 160; N = 1400; O = 1500; P = 1450; S = 1450 }
    19 = @{ D = 44417; L = "Primera"; M = 200; N = 1300; O = 1400; P = 1350; S = 1350 }
    20 = @{ D = 44260; L = "Primera"; M = 100; N = 1900; O = 2000; P = 1950; S = 1950 }
    21 = @{ D = 44350; L = "Primera"; M = 140; N = 750; O = 800; P = 775; S = 775 }
    22 = @{ D = 44351; L = "Primera"; M = 100; N = 700; O = 800; P = 750; S = 750 }
    23 = @{ D = 44351; L = "Segunda"; M = 100; N = 600; O = 700; P = 650; S = 650 }
    24 = @{ D = 44348; L = "Primera"; M = 120; N = 1000; O = 1100; P = 1050; S = 1050 }
    25 = @{ D = 44358; L = "Primera"; M = 200; N = 700; O = 800; P = 750; S = 750 }
    26 = @{ D = 44358; L = "Segunda"; M = 200; N = 600; O = 650; P = 625; S = 625 }
    27 = @{ D = 44414; L = "Primera"; M = 160; N = 1300; O = 1400; P = 1350; S = 1350 }
    28 = @{ D = 44425; L = "Primera"; M = 140; N = 1200; O = 1300; P = 1250; S = 1250 }
    29 = @{ D = 44498; L = "Segunda"; M = 100; N = 1200; O = 1300; P = 1250; S = 1250 }
    30 = @{ D = 44473; L = "Primera"; M = 160; N = 1500; O = 1600; P = 1550; S = 1550 }
    31 = @{ D = 44407; L = "Primera"; M = 200; N = 600; O = 650; P = 625; S = 625 }
}

foreach ($r in $rows.Keys) {
    $row = $rows[$r]
    $ws.Range("D$r").Value = $row.D
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
    $ws.Range("N$r").Value = $row.N
    $ws.Range("O$r").Value = $row.O
    $ws.Range("P$r").Value = $row.P
    $ws.Range("S$r").Value = $row.S
}
